$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10 and 11 (Inflammatory-Mac target-cluster rows removed)
$ws.Rows("10:11").Delete()

# Row 2
$ws.Cells.Item(2,9).Value = 0.08983953209358128
$ws.Cells.Item(2,10).Value = 0.0898395320935813
$ws.Cells.Item(2,13).Value = 2.544438666666667
$ws.Cells.Item(2,14).Value = 7.633316000000001
$ws.Cells.Item(2,15).Value = 0.201325300207035
$ws.Cells.Item(2,16).Value = 0.201325300207035
$ws.Cells.Item(2,17).Value = 0.03048491966533334
$ws.Cells.Item(2,18).Value = 0.274364276988
$ws.Cells.Item(2,19).Value = 0.01808697076919981
$ws.Cells.Item(2,20).Value = 0.01808697076919981

# Row 3
$ws.Cells.Item(3,9).Value = 0.08983953209358128
$ws.Cells.Item(3,10).Value = 0.0898395320935813
$ws.Cells.Item(3,15).Value = 0.6969390273602759
$ws.Cells.Item(3,16).Value = 0.696939027360276
$ws.Cells.Item(3,19).Value = 0.06261267611580283
$ws.Cells.Item(3,20).Value = 0.06261267611580285

# Row 4
$ws.Cells.Item(4,4).Value = 'MuSCs'
$ws.Cells.Item(4,9).Value = 0.08983953209358128
$ws.Cells.Item(4,10).Value = 0.0898395320935813
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.273916333333333
$ws.Cells.Item(4,14).Value = 3.821749
$ws.Cells.Item(4,15).Value = 0.1007969229547075
$ws.Cells.Item(4,16).Value = 0.1007969229547075
$ws.Cells.Item(4,17).Value = 0.01526279158966667
$ws.Cells.Item(4,18).Value = 0.137365124307
$ws.Cells.Item(4,19).Value = 0.009055548394723681
$ws.Cells.Item(4,20).Value = 0.009055548394723684

# Row 5
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,9).Value = 0.08983953209358128
$ws.Cells.Item(5,10).Value = 0.0898395320935813
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.01186433333333333
$ws.Cells.Item(5,14).Value = 0.035593
$ws.Cells.Item(5,15).Value = 0.0009387494779816524
$ws.Cells.Item(5,16).Value = 0.0009387494779816526
$ws.Cells.Item(5,17).Value = 0.0001421465776666667
$ws.Cells.Item(5,18).Value = 0.001279319199
$ws.Cells.Item(5,19).Value = 0.00008433681385496534
$ws.Cells.Item(5,20).Value = 0.00008433681385496537

# Row 6
$ws.Cells.Item(6,1).Value = 'MuSCs'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.121379
$ws.Cells.Item(6,8).Value = 0.364137
$ws.Cells.Item(6,9).Value = 0.9101604679064187
$ws.Cells.Item(6,10).Value = 0.9101604679064187
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.544438666666667
$ws.Cells.Item(6,14).Value = 7.633316000000001
$ws.Cells.Item(6,15).Value = 0.201325300207035
$ws.Cells.Item(6,16).Value = 0.201325300207035
$ws.Cells.Item(6,17).Value = 0.3088414209213333
$ws.Cells.Item(6,18).Value = 2.779572788292
$ws.Cells.Item(6,19).Value = 0.1832383294378352
$ws.Cells.Item(6,20).Value = 0.1832383294378352

# Row 7
$ws.Cells.Item(7,4).Value = 'FAPs'
$ws.Cells.Item(7,7).Value = 0.121379
$ws.Cells.Item(7,8).Value = 0.364137
$ws.Cells.Item(7,9).Value = 0.9101604679064187
$ws.Cells.Item(7,10).Value = 0.9101604679064187
$ws.Cells.Item(7,13).Value = 8.808225333333333
$ws.Cells.Item(7,14).Value = 26.424676
$ws.Cells.Item(7,15).Value = 0.6969390273602759
$ws.Cells.Item(7,16).Value = 0.696939027360276
$ws.Cells.Item(7,17).Value = 1.069133582734667
$ws.Cells.Item(7,18).Value = 9.622202244612
$ws.Cells.Item(7,19).Value = 0.634326351244473
$ws.Cells.Item(7,20).Value = 0.6343263512444731

# Row 8
$ws.Cells.Item(8,4).Value = 'MuSCs'
$ws.Cells.Item(8,7).Value = 0.121379
$ws.Cells.Item(8,8).Value = 0.364137
$ws.Cells.Item(8,9).Value = 0.9101604679064187
$ws.Cells.Item(8,10).Value = 0.9101604679064187
$ws.Cells.Item(8,13).Value = 1.273916333333333
$ws.Cells.Item(8,14).Value = 3.821749
$ws.Cells.Item(8,15).Value = 0.1007969229547075
$ws.Cells.Item(8,16).Value = 0.1007969229547075
$ws.Cells.Item(8,17).Value = 0.1546266906236667
$ws.Cells.Item(8,18).Value = 1.391640215613
$ws.Cells.Item(8,19).Value = 0.09174137455998378
$ws.Cells.Item(8,20).Value = 0.0917413745599838

# Row 9
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,7).Value = 0.121379
$ws.Cells.Item(9,8).Value = 0.364137
$ws.Cells.Item(9,9).Value = 0.9101604679064187
$ws.Cells.Item(9,10).Value = 0.9101604679064187
$ws.Cells.Item(9,13).Value = 0.01186433333333333
$ws.Cells.Item(9,14).Value = 0.035593
$ws.Cells.Item(9,15).Value = 0.0009387494779816524
$ws.Cells.Item(9,16).Value = 0.0009387494779816526
$ws.Cells.Item(9,17).Value = 0.001440080915666667
$ws.Cells.Item(9,18).Value = 0.012960728241
$ws.Cells.Item(9,19).Value = 0.000854412664126687
$ws.Cells.Item(9,20).Value = 0.0008544126641266873

Write-Output "Applied $($wb.Name) updates"